$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace column A raw values with "0.45 - x" formulas (x = previously
#     displayed two-decimal value), fixing the uncertainty calculation ---
$ws.Range("A2").Formula = "= 0.45 -0.1"
$ws.Range("A3").Formula = "= 0.45 -0.13"
$ws.Range("A4").Formula = "= 0.45 -0.16"
$ws.Range("A5").Formula = "= 0.45 -0.19"
$ws.Range("A6").Formula = "= 0.45 -0.22"
$ws.Range("A7").Formula = "= 0.45 -0.25"
$ws.Range("A8").Formula = "= 0.45 -0.28"
$ws.Range("A9").Formula = "= 0.45 -0.31"
$ws.Range("A10").Formula = "= 0.45 -0.34"
$ws.Range("A11").Formula = "= 0.45 -0.37"

# --- Clean up now-unused helper/plot-label columns G:H (rows 1-12) and the
#     stray style-only cell B12, clarifying the plotted-label layout ---
$ws.Range("G1:H11").Clear() | Out-Null
$ws.Range("B12").Clear() | Out-Null

# --- Update the active selection left on the sheet ---
$ws.Range("I18").Select() | Out-Null
